{"js": "// Minor changes to strategies\n//\n// 1) Remove the bullet paragraph \"Each must implement at least two\n//    different structures from the JCF.\" entirely.\n// 2) Append \" (Heuristics to change)\" to the end of the bullet\n//    \"Different strategies for the Sequence.\"\n\nconst body = context.document.body;\n\n// --- 1) Delete the \"Each must implement ...\" bullet paragraph ---\nconst removeTarget =\n  \"Each must implement at least two different structures from the JCF.\";\nconst removeResults = body.search(removeTarget, { matchCase: true });\nremoveResults.load(\"items\");\nawait context.sync();\n\nif (removeResults.items.length > 0) {\n  const para = removeResults.items[0].paragraphs.getFirst();\n  para.delete();\n  await context.sync();\n}\n\n// --- 2) Append \" (Heuristics to change)\" to the \"Different strategies\n//         for the Sequence.\" bullet ---\nconst appendTarget = \"Different strategies for the Sequence.\";\nconst appendResults = body.search(appendTarget, { matchCase: true });\nappendResults.load(\"items\");\nawait context.sync();\n\nif (appendResults.items.length > 0) {\n  const para = appendResults.items[0].paragraphs.getFirst();\n  const end = para.getRange(\"End\");\n  end.insertText(\" (Heuristics to change)\", Word.InsertLocation.end);\n  await context.sync();\n}\n", "ps1": "# Minor changes to strategies\n#\n# 1) Remove the bullet paragraph \"Each must implement at least two\n#    different structures from the JCF.\" entirely.\n# 2) Append \" (Heuristics to change)\" to the end of the bullet\n#    \"Different strategies for the Sequence.\"\n\n$d = $word.ActiveDocument\n\n# --- 1) Delete the \"Each must implement ...\" bullet paragraph ---\n$removeTarget = \"Each must implement at least two different structures from the JCF.\"\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`n\") -eq $removeTarget) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- 2) Append \" (Heuristics to change)\" to the \"Different strategies\n#         for the Sequence.\" bullet ---\n$appendTarget = \"Different strategies for the Sequence.\"\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`n\") -eq $appendTarget) {\n        $p.Range.InsertAfter(\" (Heuristics to change)\")\n        break\n    }\n}\n"}
